$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename product name (row 34, column B)
$ws.Range("B34").Value = "Best selling/recommended products"

# Status column (E) updates: "Open" -> "Deployed" or "Implemented"
$ws.Range("E6").Value = "Deployed"
$ws.Range("E9").Value = "Deployed"
$ws.Range("E13").Value = "Implemented"
$ws.Range("E31").Value = "Implemented"
$ws.Range("E32").Value = "Implemented"
$ws.Range("E33").Value = "Implemented"
$ws.Range("E34").Value = "Implemented"
$ws.Range("E36").Value = "Deployed"
$ws.Range("E38").Value = "Implemented"

# Update the selected cell in the sheet view
$ws.Range("E38").Select()
